$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 17
$ws.Range("H17").Value = 2095.1455
$ws.Range("J17").Value = 2161
$ws.Range("L17").Value = 6483
$ws.Range("N17").Value = -6819

# Row 33
$ws.Range("H33").Value = 1018.5
$ws.Range("I33").Value = 1161.0476
$ws.Range("J33").Value = 419.8
$ws.Range("K33").Value = 1161.0476
$ws.Range("L33").Value = 419.8
$ws.Range("M33").Value = -932.0476000000001
$ws.Range("N33").Value = -877.8

# Row 74
$ws.Range("H74").Value = 4248.04
$ws.Range("I74").Value = 4588.5557
$ws.Range("J74").Value = 4056.5
$ws.Range("K74").Value = 4588.5557
$ws.Range("L74").Value = 4056.5
$ws.Range("M74").Value = -3652.5557
$ws.Range("N74").Value = -5928.5

# Row 77
$ws.Range("H77").Value = 4248.04
$ws.Range("I77").Value = 4588.5557
$ws.Range("J77").Value = 4056.5
$ws.Range("K77").Value = 22942.7785
$ws.Range("L77").Value = 20282.5
$ws.Range("M77").Value = -18262.7785
$ws.Range("N77").Value = -29642.5

# Row 112
$ws.Range("H112").Value = 1391.1305
$ws.Range("I112").Value = 500
$ws.Range("J112").Value = 1578.7368
$ws.Range("K112").Value = 1500
$ws.Range("L112").Value = 4736.2104
$ws.Range("M112").Value = -392
$ws.Range("N112").Value = -6952.2104

# Row 127
$ws.Range("H127").Value = 837.4737
$ws.Range("I127").Value = 536.3333
$ws.Range("J127").Value = 1108.5
$ws.Range("K127").Value = 1608.9999
$ws.Range("L127").Value = 3325.5
$ws.Range("M127").Value = 3351.0001
$ws.Range("N127").Value = -13245.5

# Row 132
$ws.Range("H132").Value = 5876
$ws.Range("I132").Value = 6922.7896
$ws.Range("J132").Value = 4346.077
$ws.Range("K132").Value = 20768.3688
$ws.Range("L132").Value = 13038.231
$ws.Range("M132").Value = -18238.3688
$ws.Range("N132").Value = -18098.231

# Row 137
$ws.Range("H137").Value = 1238.6719
$ws.Range("I137").Value = 830.72095
$ws.Range("K137").Value = 2492.16285
$ws.Range("M137").Value = 57.83714999999984

$ws = $wb.Worksheets.Item("ARM")
# Row 32
$ws.Range("H32").Value = 10190.925
$ws.Range("I32").Value = 3172.25
$ws.Range("J32").Value = 26567.834
$ws.Range("K32").Value = 3172.25
$ws.Range("L32").Value = 26567.834
$ws.Range("M32").Value = -2885.25
$ws.Range("N32").Value = -27141.834

# Row 74
$ws.Range("H74").Value = 7145608.5
$ws.Range("I74").Value = 8930877
$ws.Range("K74").Value = 8930877
$ws.Range("M74").Value = -8930003

# Row 77
$ws.Range("H77").Value = 7145608.5
$ws.Range("I77").Value = 8930877
$ws.Range("K77").Value = 44654385
$ws.Range("M77").Value = -44650017

# Row 110
$ws.Range("H110").Value = 5957.815
$ws.Range("I110").Value = 7024.6
$ws.Range("J110").Value = 2909.8572
$ws.Range("K110").Value = 7024.6
$ws.Range("L110").Value = 2909.8572
$ws.Range("M110").Value = -4979.6
$ws.Range("N110").Value = -6999.8572

# Row 132
$ws.Range("H132").Value = 2769.889
$ws.Range("I132").Value = 2652.5
$ws.Range("J132").Value = 3004.6667
$ws.Range("K132").Value = 7957.5
$ws.Range("L132").Value = 9014.000100000001
$ws.Range("M132").Value = -5427.5
$ws.Range("N132").Value = -14074.0001

$ws = $wb.Worksheets.Item("BSM")
# Row 86
$ws.Range("H86").Value = 11766465
$ws.Range("I86").Value = 15386101
$ws.Range("J86").Value = 2649.25
$ws.Range("K86").Value = 15386101
$ws.Range("L86").Value = 2649.25
$ws.Range("M86").Value = -15384978
$ws.Range("N86").Value = -4895.25

# Row 89
$ws.Range("H89").Value = 11766465
$ws.Range("I89").Value = 15386101
$ws.Range("J89").Value = 2649.25
$ws.Range("K89").Value = 76930505
$ws.Range("L89").Value = 13246.25
$ws.Range("M89").Value = -76924889
$ws.Range("N89").Value = -24478.25

# Row 107
$ws.Range("H107").Value = 876.0263
$ws.Range("I107").Value = 628.28
$ws.Range("K107").Value = 628.28
$ws.Range("M107").Value = 1291.72

$ws = $wb.Worksheets.Item("CRP")
# Row 31
$ws.Range("H31").Value = 4371066
$ws.Range("I31").Value = 7731732
$ws.Range("J31").Value = 2199.45
$ws.Range("K31").Value = 7731732
$ws.Range("L31").Value = 2199.45
$ws.Range("M31").Value = -7731437
$ws.Range("N31").Value = -2789.45

# Row 34
$ws.Range("H34").Value = 4371066
$ws.Range("I34").Value = 7731732
$ws.Range("J34").Value = 2199.45
$ws.Range("K34").Value = 7731732
$ws.Range("L34").Value = 2199.45
$ws.Range("M34").Value = -7731530
$ws.Range("N34").Value = -2603.45

# Row 132
$ws.Range("H132").Value = 1485.1111
$ws.Range("I132").Value = 1252.7059
$ws.Range("J132").Value = 1880.2
$ws.Range("K132").Value = 3758.1177
$ws.Range("L132").Value = 5640.6
$ws.Range("M132").Value = -1228.1177
$ws.Range("N132").Value = -10700.6

$ws = $wb.Worksheets.Item("CUL")
# Row 121
$ws.Range("H121").Value = 6894.029
$ws.Range("I121").Value = 6955.2
$ws.Range("J121").Value = 6848.15
$ws.Range("K121").Value = 20865.6
$ws.Range("L121").Value = 20544.45
$ws.Range("M121").Value = -19555.6
$ws.Range("N121").Value = -23164.45

$ws = $wb.Worksheets.Item("LTW")
# Row 82
$ws.Range("H82").Value = 1705.5
$ws.Range("I82").Value = 1400
$ws.Range("J82").Value = 1875.2222
$ws.Range("K82").Value = 1400
$ws.Range("L82").Value = 1875.2222
$ws.Range("M82").Value = -1039
$ws.Range("N82").Value = -2597.2222

# Row 85
$ws.Range("H85").Value = 1705.5
$ws.Range("I85").Value = 1400
$ws.Range("J85").Value = 1875.2222
$ws.Range("K85").Value = 1400
$ws.Range("L85").Value = 1875.2222
$ws.Range("M85").Value = -152
$ws.Range("N85").Value = -4371.2222

# Row 132
$ws.Range("H132").Value = 22331684
$ws.Range("I132").Value = 41682920
$ws.Range("J132").Value = 3335
$ws.Range("K132").Value = 125048760
$ws.Range("L132").Value = 10005
$ws.Range("M132").Value = -125046230
$ws.Range("N132").Value = -15065

$ws = $wb.Worksheets.Item("WVR")
# Row 81
$ws.Range("H81").Value = 41668580
$ws.Range("I81").Value = 111113300
$ws.Range("J81").Value = 1750.5333
$ws.Range("K81").Value = 222226600
$ws.Range("L81").Value = 3501.0666
$ws.Range("M81").Value = -222225539
$ws.Range("N81").Value = -5623.0666

# Row 84
$ws.Range("H84").Value = 41668580
$ws.Range("I84").Value = 111113300
$ws.Range("J84").Value = 1750.5333
$ws.Range("K84").Value = 1111133000
$ws.Range("L84").Value = 17505.333
$ws.Range("M84").Value = -1111127696
$ws.Range("N84").Value = -28113.333

